# Hornerstown S9000 - config generator / tkinter test script update
# Re-labels the point list (columns A/C/D) and shifts the value-type
# column (old F -> new E, old D -> new D) to make room for new MW/MVAR/
# amps points appended in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1-10: point name (A), open/close or off/on (D, was F), and the
#     new MW/MVAR/amps point (E, was D) ---------------------------------
$ws.Range("A1").Value  = "BK 1 VCB"
$ws.Range("D1").Value  = "OFF"
$ws.Range("E1").Value  = "SPARE"
$ws.Range("F1").ClearContents()

$ws.Range("A2").Value  = "47417 VCB"
$ws.Range("D2").Value  = "ON"
$ws.Range("E2").Value  = "SPARE"
$ws.Range("F2").ClearContents()

$ws.Range("A3").Value  = "BT VCB"
$ws.Range("D3").Value  = "OPEN"
$ws.Range("E3").Value  = "BK 1 MW"
$ws.Range("F3").ClearContents()

$ws.Range("A4").Value  = "BK 2 VCB"
$ws.Range("D4").Value  = "CLOSE"
$ws.Range("E4").Value  = "BK 1 MVAR"
$ws.Range("F4").ClearContents()

$ws.Range("A5").Value  = "47418 VCB"
$ws.Range("D5").Value  = "OPEN"
$ws.Range("E5").Value  = "BK 2 MW"
$ws.Range("F5").ClearContents()

$ws.Range("A6").Value  = "47419 VCB"
$ws.Range("D6").Value  = "CLOSE"
$ws.Range("E6").Value  = "BK 2 MVAR"
$ws.Range("F6").ClearContents()

$ws.Range("A7").Value  = "47416 VCB"
$ws.Range("D7").Value  = "OPEN"
$ws.Range("E7").Value  = "47416 phase A amps"
$ws.Range("F7").ClearContents()

$ws.Range("A8").Value  = "47415 VCB"
$ws.Range("D8").Value  = "CLOSE"
$ws.Range("E8").Value  = "47416 phase B amps"
$ws.Range("F8").ClearContents()

$ws.Range("A9").Value  = "DIAL IN ACCESS"
$ws.Range("D9").Value  = "OPEN"
$ws.Range("E9").Value  = "47416 phase C amps"
$ws.Range("F9").ClearContents()

$ws.Range("A10").Value = "BK 1 CKT INT CI-2"
$ws.Range("D10").Value = "CLOSE"
$ws.Range("E10").Value = "47417 phase A amps"
$ws.Range("F10").ClearContents()

# --- Rows 11-14: additional new points only feed column E, A/B/C unchanged
$ws.Range("E11").Value = "47417 phase B amps"
$ws.Range("E12").Value = "47417 phase C amps"
$ws.Range("A14").Value = "SPARE "

# --- Rows 15-27: relabeled point names in column A ----------------------
$ws.Range("A15").Value = "I87 LBSW"
$ws.Range("A16").Value = "I87 LBSW"
$ws.Range("A17").Value = "NO 6 LBSW"
$ws.Range("A18").Value = "NO 6 LBSW"
$ws.Range("A19").Value = "A53 LBSW"
$ws.Range("A20").Value = "A53 LBSW"
$ws.Range("A21").Value = "BK 1 ALRM"
$ws.Range("A22").Value = "BK 2 ALRM"
$ws.Range("A23").Value = "CAP 1 VCB"
$ws.Range("A24").Value = "DX RELAY FAIL"
$ws.Range("A25").Value = "LOSS OF POTENTIAL "
$ws.Range("A26").Value = "BK 1 CKT INT CI-2/SEL TRBL  "
$ws.Range("A27").Value = "ADAPTIVE RELAYING"

# --- Rows 29-42: trailing "SPARE" placeholders renamed to "UNDEFINED" ---
$ws.Range("A29").Value = "UNDEFINED"
$ws.Range("A30").Value = "UNDEFINED"
$ws.Range("A31").Value = "UNDEFINED"
$ws.Range("A32").Value = "UNDEFINED"
$ws.Range("A33").Value = "UNDEFINED"
$ws.Range("A34").Value = "UNDEFINED"
$ws.Range("A35").Value = "UNDEFINED"
$ws.Range("A36").Value = "UNDEFINED"
$ws.Range("A37").Value = "UNDEFINED"
$ws.Range("A38").Value = "UNDEFINED"
$ws.Range("A39").Value = "UNDEFINED"
$ws.Range("A40").Value = "UNDEFINED"
$ws.Range("A41").Value = "UNDEFINED"
$ws.Range("A42").Value = "UNDEFINED"

# --- Formatting: column E now carries data too, so widen it like D, and
#     drop the leftover number-format style that used to sit (unused) on
#     E1 so it reverts to the sheet's normal style -----------------------
$ws.Range("E1").Style = "Normal"
$ws.Columns.Item(5).ColumnWidth = 18

# --- Selection / scroll position reset (matches the saved view) ---------
$ws.Range("A15").Select()
